$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are forced to text via a temporary "@" (Text)
# number format so Excel does not auto-coerce numeric-looking strings like
# "593.89" into floating point numbers; the style is then reset back to
# "Normal" so no residual formatting diff is left on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.048.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.034.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.28%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.89%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.030.56"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +17.32%  "
$ws.Range("E11").Value = "  +4.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.463"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.78%  "
$ws.Range("E13").Value = "  +4.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.60%  "
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.534.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.982.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.031.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "453.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.699"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +12.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.10%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +13.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.95%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.111"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0864"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.59%  "
$ws.Range("E36").Value = "  +3.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.24%  "
$ws.Range("E39").Value = "  +9.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("E42").Value = "  +2.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.312"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +18.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +16.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "396.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0360"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.721.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.84%  "
$ws.Range("E51").Value = "  +8.57%  "
